$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the latest GitHub Actions scrape.
# Price column (D) values are forced to text (NumberFormat "@") so values such
# as "594.79" are stored as strings (matching the source inlineStr cells)
# instead of being auto-converted to numbers by Excel; the style is then
# reset to "Normal" so no extra formatting is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.776.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.075.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.46%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.076.24'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("E10").Value = '  -0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.89'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("E12").Value = '  -1.75%  '
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.97'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.68%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000238'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.02%  '
$ws.Range("E15").Value = '  +1.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.584.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.743.67'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.073.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '489.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.69%  '
$ws.Range("E22").Value = '  -1.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.32%  '
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.90%  '
$ws.Range("E27").Value = '  +7.54%  '
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.66%  '
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.13%  '
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.36'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("E34").Value = '  -2.69%  '
$ws.Range("E35").Value = '  +1.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0824'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.36%  '
$ws.Range("E37").Value = '  -1.80%  '
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.29'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.03%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.23'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.26'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '439.36'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.292'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.114'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.22%  '
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.843.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.22'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.02%  '
